$d = $word.ActiveDocument

$replacements = @(
    @("2023-11-23 Thursday", "2023-11-24 Friday"),
    @("89÷9=", "58÷3="),
    @("32÷3=", "96÷3="),
    @("37÷7=", "92÷3="),
    @("45÷4=", "58÷4="),
    @("53÷6=", "21÷5="),
    @("76÷4=", "10÷2="),
    @("80÷5=", "71÷2="),
    @("75÷2=", "64÷2="),
    @("52÷9=", "39÷2="),
    @("44÷4=", "73÷4="),
    @("96÷2=", "80÷3="),
    @("92÷4=", "57÷2="),
    @("71÷6=", "37÷6="),
    @("38÷7=", "61÷3="),
    @("39÷3=", "93÷5="),
    @("43÷7=", "75÷3="),
    @("71÷5=", "98÷5="),
    @("13÷9=", "67÷9="),
    @("76÷5=", "18÷7="),
    @("30÷4=", "77÷5="),
    @("94÷8=", "50÷6="),
    @("21÷9=", "49÷6="),
    @("57÷4=", "51÷2="),
    @("94÷5=", "40÷8="),
    @("92÷5=", "13÷3="),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

